$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column widths for F and G
$ws.Columns.Item(6).ColumnWidth = 15.88671875
$ws.Columns.Item(7).ColumnWidth = 13.44140625

# Header row
$ws.Range("F1").Value = "Форма контроля"
$ws.Range("G1").Value = "Доп. Контроль"
$ws.Range("F1:G1").Font.Bold = $true

# Section 1 (rows 3-6)
$ws.Range("F3").Value = "Зачет"
$ws.Range("G3").Value = "Нет"

$ws.Range("F4").Value = "Экз."
$ws.Range("G4").Value = "Нет"

$ws.Range("F5").Value = "Зачет"
$ws.Range("G5").Value = "Курсовая"

$ws.Range("F6").Value = "Экз."
$ws.Range("G6").Value = "Нет"

# Section 2 (rows 9-13)
$ws.Range("F9").Value = "Экз."
$ws.Range("G9").Value = "Проект"

$ws.Range("F10").Value = "Экз."
$ws.Range("G10").Value = "Нет"

$ws.Range("F11").Value = "Экз."
$ws.Range("G11").Value = "Нет"

$ws.Range("D12").Value = 16
$ws.Range("F12").Value = "Зачет"
$ws.Range("G12").Value = "Курсовая"

$ws.Range("D13").Value = 16
$ws.Range("F13").Value = "Зачет"
$ws.Range("G13").Value = "Курсовая"

# Section 3 (rows 16-23)
$ws.Range("F16").Value = "Экз."
$ws.Range("G16").Value = "Нет"

$ws.Range("F17").Value = "Зачет"
$ws.Range("G17").Value = "Проект"

$ws.Range("F18").Value = "Зачет"
$ws.Range("G18").Value = "Нет"

$ws.Range("F19").Value = "Зачет"
$ws.Range("G19").Value = "Нет"

$ws.Range("D20").Value = 16
$ws.Range("F20").Value = "Экз."
$ws.Range("G20").Value = "Курсовая"

$ws.Range("F21").Value = "Зачет"
$ws.Range("G21").Value = "Нет"

$ws.Range("B22").Value = "DevOPS"
$ws.Range("F22").Value = "Зачет"
$ws.Range("G22").Value = "Нет"

$ws.Range("D23").Value = 16
$ws.Range("F23").Value = "Экз."
$ws.Range("G23").Value = "Проект"

# New row 25: Факультативы section header, merged A25:E25
$ws.Range("A25").Value = "Факультативы"
$ws.Range("A25:E25").Merge()
$ws.Range("A25:E25").HorizontalAlignment = -4108

# Update selection to E14
$ws.Range("E14").Select()
